$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data row to append: date serial 46008 (12/17/2025), error count 9
$newRow = 36
$prevRow = $newRow - 1

# Copy the previous date cell's format (numFmtId 14, mm/dd/yyyy) onto the new cell
$ws.Range("A" + $prevRow).Copy()
$ws.Range("A" + $newRow).PasteSpecial(-4122)

$ws.Cells.Item($newRow, 1).Value = 46008
$ws.Cells.Item($newRow, 2).Value = 9

# Update the active selection to mirror the appended row, like Excel would after data entry
$ws.Range("A36:B36").Select()
